# Update the BLE tracking rows in both sheets: Floodplain status (E) for the
# rows that moved to Approved / In Backcheck, and convert the P01 GDB (G),
# RAW Grids (H) and Model Complete (N) date-as-text columns into real Excel
# date serials with a date-time number format.
$wb = $excel.ActiveWorkbook

$rowUpdates = @(
    @{ Row=2; E="Approved"; G=45350; H=45386; N=45301 }
    @{ Row=3; E="Approved"; G=45290; H=45366; N=45241 }
    @{ Row=4; E="Approved"; G=45434; H=45397; N=45241 }
    @{ Row=5; E="Approved"; G=45434; H=45429; N=45272 }
    @{ Row=6; E="Approved"; G=45434; H=45421; N=45301 }
    @{ Row=7; E="Approved"; G=45260; H=45345; N=45211 }
    @{ Row=8; E="Approved"; G=45321; H=45366; N=45272 }
    @{ Row=9; E="Approved"; G=45486; H=45457; N=45443 }
    @{ Row=10; E="In Backcheck"; G=45381; H=45402; N=45332 }
    @{ Row=11; E="Approved"; G=45260; H=45345; N=45211 }
    @{ Row=12; E="Approved"; G=45290; H=45397; N=45241 }
    @{ Row=13; G=45290; H=45375; N=45241 }
    @{ Row=14; G=45508; H=45508; N=45498 }
    @{ Row=15; G=45348; H=45335; N=45150 }
    @{ Row=16; G=45137; H=45169; N=45088 }
    @{ Row=17; G=45446; H=45356; N=45241 }
    @{ Row=18; G=45350; H=45377; N=45301 }
    @{ Row=19; G=45381; H=45436; N=45332 }
    @{ Row=20; G=45503; H=45503; N=45494 }
    @{ Row=21; G=45484; H=45455; N=45481 }
    @{ Row=22; G=45544; H=45485; N=45481 }
    @{ Row=23; G=45505; H=45505; N=45363 }
    @{ Row=24; G=45441; H=45440; N=45332 }
    @{ Row=25; G=45509; H=45510; N=45424 }
    @{ Row=26; G=45509; H=45511; N=45489 }
    @{ Row=27; G=45524; H=45510; N=45505 }
)

$dateColumns = @("G", "H", "N")
$dateFormat = "YYYY-MM-DD HH:MM:SS"
$numFmtSeeded = $false

foreach ($sheetName in @("Tracking_Main_values", "Tracking_Main")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($item in $rowUpdates) {
        if ($item.ContainsKey("E")) {
            $ws.Range("E" + $item.Row).Value = $item.E
        }

        foreach ($col in $dateColumns) {
            $cell = $ws.Range($col + $item.Row)
            $cell.Value = $item[$col]
            if (-not $numFmtSeeded) {
                # Seed numFmtId 164 with a lowercase code (unused by any cell,
                # mirrors the original edit) so the real format below lands on 165.
                $cell.NumberFormat = "yyyy-mm-dd h:mm:ss"
                $numFmtSeeded = $true
            }
            $cell.NumberFormat = $dateFormat
        }
    }
}
